$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the LR-pair data table from 4 rows (2-5) to 9 rows (2-10),
# adding "ECs" as a third sending/target cluster alongside FAPs and sCs,
# and updating all recomputed expression/specificity statistics.
$data = @(
    @("FAPs", "Vegfc", "Vipr2", "FAPs", 2.0, 0.6666666666666666, 3.911257666666666, 11.733773, 0.4115343446855154, 0.4115343446855154, 3.0, 1.0, 0.5692973333333333, 1.707892, 0.07854540531141245, 0.07854540531141245, 2.226668559612889, 20.040017036516, 0.03232413190289032, 0.03232413190289032),
    @("FAPs", "Vegfc", "Vipr2", "ECs", 2.0, 0.6666666666666666, 3.911257666666666, 11.733773, 0.4115343446855154, 0.4115343446855154, 3.0, 1.0, 3.384345, 10.153035, 0.4669348232885666, 0.4669348232885665, 13.237045327895, 119.133407951055, 0.1921597165129072, 0.1921597165129071),
    @("FAPs", "Vegfc", "Vipr2", "sCs", 2.0, 0.6666666666666666, 3.911257666666666, 11.733773, 0.4115343446855154, 0.4115343446855154, 2.0, 0.6666666666666666, 3.294360666666666, 9.883082, 0.4545197714000211, 0.4545197714000211, 12.88509341426511, 115.965840728386, 0.1870504962697179, 0.1870504962697179),
    @("ECs", "Vegfc", "Vipr2", "FAPs", 3.0, 1.0, 4.100181333333333, 12.300544, 0.4314124974392592, 0.4314124974392592, 3.0, 1.0, 0.5692973333333333, 1.707892, 0.07854540531141245, 0.07854540531141245, 2.334222299249778, 21.008000693248, 0.0338854694677753, 0.0338854694677753),
    @("ECs", "Vegfc", "Vipr2", "ECs", 3.0, 1.0, 4.100181333333333, 12.300544, 0.4314124974392592, 0.4314124974392592, 3.0, 1.0, 3.384345, 10.153035, 0.4669348232885666, 0.4669348232885665, 13.87642819456, 124.88785375104, 0.2014415182562797, 0.2014415182562797),
    @("ECs", "Vegfc", "Vipr2", "sCs", 3.0, 1.0, 4.100181333333333, 12.300544, 0.4314124974392592, 0.4314124974392592, 2.0, 0.6666666666666666, 3.294360666666666, 9.883082, 0.4545197714000211, 0.4545197714000211, 13.50747611073422, 121.567284996608, 0.1960855097152043, 0.1960855097152043),
    @("sCs", "Vegfc", "Vipr2", "FAPs", 3.0, 1.0, 1.492646666666667, 4.47794, 0.1570531578752254, 0.1570531578752254, 3.0, 1.0, 0.5692973333333333, 1.707892, 0.07854540531141245, 0.07854540531141245, 0.8497597669422222, 7.64783790248, 0.01233580394074683, 0.01233580394074683),
    @("sCs", "Vegfc", "Vipr2", "ECs", 3.0, 1.0, 1.492646666666667, 4.47794, 0.1570531578752254, 0.1570531578752254, 3.0, 1.0, 3.384345, 10.153035, 0.4669348232885666, 0.4669348232885665, 5.0516312831, 45.4646815479, 0.07333358851937971, 0.07333358851937971),
    @("sCs", "Vegfc", "Vipr2", "sCs", 3.0, 1.0, 1.492646666666667, 4.47794, 0.1570531578752254, 0.1570531578752254, 2.0, 0.6666666666666666, 3.294360666666666, 9.883082, 0.4545197714000211, 0.4545197714000211, 4.917316467897778, 44.25584821108, 0.07138376541509886, 0.07138376541509886)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $rowValues[$j]
    }
}
